$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing data rows 2-5, pushing the
# current rows 2-5 down to become rows 6-9.
$ws.Range("A2:A5").EntireRow.Insert()

# The newly inserted rows picked up the header row's formatting
# (bold, centered). Reset to the default (unstyled) formatting that
# the target data rows use.
$ws.Range("A2:T5").ClearFormats()

# Restore the date-formatted style (custom date number format) onto
# column D for the new rows, copying it from the (now shifted) row 6
# which still carries the original date style.
$ws.Range("D6").Copy()
$ws.Range("D2:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: Especial
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44908
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104004
$ws.Range("J2").Value = "Níspero"
$ws.Range("K2").Value = "Golden Nugget"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 110
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1400
$ws.Range("T2").Value = 5

# Row 3: Primera
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44908
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104004
$ws.Range("J3").Value = "Níspero"
$ws.Range("K3").Value = "Golden Nugget"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1200
$ws.Range("T3").Value = 5

# Row 4: Segunda
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44908
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104004
$ws.Range("J4").Value = "Níspero"
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 5

# Row 5: Tercera
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44908
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104004
$ws.Range("J5").Value = "Níspero"
$ws.Range("K5").Value = "Golden Nugget"
$ws.Range("L5").Value = "Tercera"
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 4000
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 4000
$ws.Range("Q5").Value = "$/bandeja 5 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 800
$ws.Range("T5").Value = 5
